# Add a new "13-sep" column (BS) to the right of the existing "12-sep"
# column (BR) on Sheet1, carrying one value per data row (2-18), then move
# the active selection one column to the right (BT6 -> BU6), matching how
# Excel shifts the cursor after data is typed into a freshly-used column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell for the new date column.
$ws.Range("BS1").Value = "13-sep"

# New column's data values, one per row.
$bsValues = @{
    2  = 0
    3  = 17.593543315544416
    4  = 16.438023609530195
    5  = 15.554346852534602
    6  = 0
    7  = 11.342060131097806
    8  = 15.277172696949297
    9  = 14.954685055067683
    10 = 16.303061899648281
    11 = 15.875530711263007
    12 = 0
    13 = 10.635555535039101
    14 = 0
    15 = 0
    16 = 8.8127046186131821
    17 = 0
    18 = 0
}

foreach ($row in $bsValues.Keys) {
    $ws.Cells.Item($row, 71).Value = $bsValues[$row]
}

# Move the selection the same way Excel does once the new column is filled.
$ws.Range("BU6").Select()
